$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Archetypes")

# --- Insert a new "code" column before the existing description column (D) ---
# This shifts the old column D (description) to E and the old column H to I,
# and keeps all same-row relative formula references (e.g. C2) intact.
$ws.Columns.Item(4).Insert()

# --- Header row ---
$ws.Range("D1").Value = "code"

# --- Row 2: windpark-small-offshore ---
$ws.Range("A2").Value = "windpark-small-offshore"
$ws.Range("B2").Formula = "='Offshore wind'!I24*C2"
$ws.Range("C2").Value = 20
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "small offshore wind park"

# --- Row 3: solarpark-small ---
$ws.Range("A3").Value = "solarpark-small"
$ws.Range("B3").Formula = "='Solar Photovoltaic'!B41*Archetypes!C3"
$ws.Range("B3").Style = "Comma"
$ws.Range("C3").Value = 15
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "small solar photovoltaic power station"

# --- Row 4: windpark-small-onshore ---
$ws.Range("A4").Value = "windpark-small-onshore"
$ws.Range("B4").Formula = "='Onshore wind'!B34*Archetypes!C4"
$ws.Range("C4").Value = 5
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "small onshore wind park"

# --- Row 5: windpark-medium-offshore (keeps the formatted-but-empty I5 cell, style preserved from old H4) ---
$ws.Range("A5").Value = "windpark-medium-offshore"
$ws.Range("B5").Formula = "='Offshore wind'!I24*C5"
$ws.Range("C5").Value = 125
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = "medium offshore wind park"

# --- Row 6: solarpark-medium ---
$ws.Range("A6").Value = "solarpark-medium"
$ws.Range("B6").Formula = "='Solar Photovoltaic'!B41*Archetypes!C6"
$ws.Range("C6").Value = 29
$ws.Range("D6").Value = 2
$ws.Range("E6").Value = "medium solar photovoltaic power station"

# --- Row 7: windpark-medium-onshore ---
$ws.Range("A7").Value = "windpark-medium-onshore"
$ws.Range("B7").Formula = "='Onshore wind'!B34*Archetypes!C7"
$ws.Range("C7").Value = 9
$ws.Range("D7").Value = 2
$ws.Range("E7").Value = "medium onshore wind park"

# --- Row 8: windpark-large-offshore ---
$ws.Range("A8").Value = "windpark-large-offshore"
$ws.Range("B8").Formula = "='Offshore wind'!I24*C8"
$ws.Range("C8").Value = 120
$ws.Range("D8").Value = 3
$ws.Range("E8").Value = "large offshore wind park"

# --- Row 9: solarpark-large ---
$ws.Range("A9").Value = "solarpark-large"
$ws.Range("B9").Formula = "='Solar Photovoltaic'!B41*Archetypes!C9"
$ws.Range("C9").Value = 47
$ws.Range("D9").Value = 3
$ws.Range("E9").Value = "large solar photovoltaic power station"

# --- Row 10: windpark-large-onshore ---
$ws.Range("A10").Value = "windpark-large-onshore"
$ws.Range("B10").Formula = "='Onshore wind'!B34*Archetypes!C10"
$ws.Range("C10").Value = 19
$ws.Range("D10").Value = 3
$ws.Range("E10").Value = "large onshore wind park"

# --- Row 11: solar-generating-district ---
$ws.Range("A11").Value = "solar-generating-district"
$ws.Range("B11").Value = 0
$ws.Range("C11").Value = 0.25
$ws.Range("D11").Value = 4
$ws.Range("E11").Value = "solar power generating district"

# --- Column D width to match column C ---
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Selection as last left by the editor ---
$ws.Range("F7").Select()
